# Updated cryptos list on Mon Nov 20 08:47:24 UTC 2023 with GitHub Actions
#
# Applies the new Price (column D) and Volume(1h) (column E) values scraped
# for each coin, and swaps rows 13/14 (Avalanche <-> WrappedliquidstakedEther2.0)
# whose relative ranking changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells are plain decimal-looking strings (e.g. "1.39",
# "247.10", "0.0806") that Excel's COM layer would otherwise auto-convert to
# numeric values (losing trailing zeros / text formatting) when assigned via
# .Value. Forcing a Text number format before the assignment - and then
# resetting the cell style back to Normal afterwards so no stray style index
# is left on the cell - keeps them stored as literal text, matching the
# original workbook's inline-string cells.
$priceCells = @(
    "D2","D3","D5","D6","D7","D9","D10","D12","D13","D14","D15","D16","D17",
    "D18","D19","D20","D21","D22","D25","D27","D28","D29","D30","D31","D32",
    "D35","D36","D39","D40","D43","D44","D45","D46","D48","D49","D50"
)
foreach ($c in $priceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "37.148.77"
$ws.Range("D3").Value = "2.013.84"
$ws.Range("D5").Value = "247.10"
$ws.Range("D6").Value = "0.625"
$ws.Range("D7").Value = "59.81"
$ws.Range("D9").Value = "0.389"
$ws.Range("D10").Value = "0.0806"
$ws.Range("D12").Value = "15.11"

# Rows 13 and 14 swap their coin content (WrappedliquidstakedEther2.0 moved
# above Avalanche in the ranking) - update every column for both rows.
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.314.69"
$ws.Range("E13").Value = "  +3.58%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "22.29"
$ws.Range("E14").Value = "  +1.56%  "

$ws.Range("D15").Value = "0.849"
$ws.Range("D16").Value = "5.52"
$ws.Range("D17").Value = "2.023.59"
$ws.Range("D18").Value = "37.051.15"
$ws.Range("D19").Value = "70.46"
$ws.Range("D20").Value = "0.0₃0868"
$ws.Range("D21").Value = "5.22"
$ws.Range("D22").Value = "230.52"
$ws.Range("D25").Value = "2.36"
$ws.Range("D27").Value = "163.84"
$ws.Range("D28").Value = "0.137"
$ws.Range("D29").Value = "19.72"
$ws.Range("D30").Value = "1.39"
$ws.Range("D31").Value = "0.121"
$ws.Range("D32").Value = "4.81"
$ws.Range("D35").Value = "2.44"
$ws.Range("D36").Value = "3.51"
$ws.Range("D39").Value = "5.40"
$ws.Range("D40").Value = "0.0987"
$ws.Range("D43").Value = "0.0214"
$ws.Range("D44").Value = "16.64"
$ws.Range("D45").Value = "92.03"
$ws.Range("D46").Value = "1.384.64"
$ws.Range("D48").Value = "7.51"
$ws.Range("D49").Value = "2.11"
$ws.Range("D50").Value = "47.08"

foreach ($c in $priceCells) {
    $ws.Range($c).Style = "Normal"
}

# Column E ("Volume(1h)") percentage strings already contain non-numeric
# characters (the leading/trailing spaces and '%' sign), so Excel keeps them
# as text without any extra handling.
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("E3").Value = "  +2.69%  "
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("E6").Value = "  +1.21%  "
$ws.Range("E7").Value = "  -2.97%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +2.67%  "
$ws.Range("E10").Value = "  +1.20%  "
$ws.Range("E11").Value = "  +1.48%  "
$ws.Range("E12").Value = "  +5.95%  "
$ws.Range("E15").Value = "  +1.38%  "
$ws.Range("E16").Value = "  +4.06%  "
$ws.Range("E17").Value = "  +3.22%  "
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("E21").Value = "  +2.63%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +1.50%  "
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("E26").Value = "  +2.38%  "
$ws.Range("E27").Value = "  +2.02%  "
$ws.Range("E28").Value = "  -3.31%  "
$ws.Range("E29").Value = "  +1.37%  "
$ws.Range("E30").Value = "  +6.34%  "
$ws.Range("E31").Value = "  +1.42%  "
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("E33").Value = "  +6.69%  "
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("E35").Value = "  +7.21%  "
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("E38").Value = "  +1.88%  "
$ws.Range("E39").Value = "  -3.14%  "
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("E41").Value = "  +0.64%  "
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("E44").Value = "  +3.28%  "
$ws.Range("E45").Value = "  +3.50%  "
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("E47").Value = "  +2.17%  "
$ws.Range("E48").Value = "  +4.85%  "
$ws.Range("E49").Value = "  +13.72%  "
$ws.Range("E50").Value = "  +2.22%  "
$ws.Range("E51").Value = "  +0.24%  "
